$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Total Load Profile for All Users for: 2023-01-01 to 2023-01-20"

$ws.Range("B3").Value = 0.5786936985157896
$ws.Range("B4").Value = 0.5553883948421053
$ws.Range("B5").Value = 0.509553242368421
$ws.Range("B6").Value = 0.3973893431
$ws.Range("B7").Value = 0.2675126149052631
$ws.Range("B8").Value = 0.317891004336842
$ws.Range("B9").Value = 0.3430626234894737
$ws.Range("B10").Value = 0.3826667084157895
$ws.Range("B11").Value = 0.439236928268421
$ws.Range("B12").Value = 0.4851446510263158
$ws.Range("B13").Value = 0.5334119425789475
$ws.Range("B14").Value = 0.5476421052631579
$ws.Range("B15").Value = 0.5462532337631579
$ws.Range("B16").Value = 0.5430888715000001
$ws.Range("B17").Value = 0.5285
$ws.Range("B18").Value = 0.5390162064263158
$ws.Range("B19").Value = 0.8446761897578948
$ws.Range("B20").Value = 1.428286588378948
$ws.Range("B21").Value = 1.437801871952632
$ws.Range("B22").Value = 1.217591516178947
$ws.Range("B23").Value = 0.8926623312947368
$ws.Range("B24").Value = 0.7417273265052632
$ws.Range("B25").Value = 0.6570347317736843
$ws.Range("B26").Value = 0.6038154944736842
